# Updates cryptos list with refreshed prices and 1h volume percentages.
# Rows 31/32 and 45/46 also swap which coin occupies each rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.043.92"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "2.648.08"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "530.77"
$ws.Range("E5").Value = "  +4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.82"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "3.104.88"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").Value = "61.040.28"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.92"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("D17").Value = "2.643.29"
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.77"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "354.75"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.69"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.73"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.431"
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").Value = "0.0₃0867"
$ws.Range("E27").Value = "  +3.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  +7.25%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.52"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.63"
$ws.Range("E32").Value = "  +4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.59"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.17"
$ws.Range("E34").Value = "  +4.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.21"
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.923"
$ws.Range("E36").Value = "  +9.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.894"
$ws.Range("E37").Value = "  +1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "308.56"
$ws.Range("E38").Value = "  +4.70%  "
$ws.Range("E39").Value = "  +1.55%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.646"
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.05"
$ws.Range("E45").Value = "  +3.57%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.86"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.37"
$ws.Range("E48").Value = "  +8.59%  "
$ws.Range("D50").Value = "1.990.10"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +2.59%  "

Write-Host "Updated cryptos list"
